$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the new "EXTRA" block (row 21) ---
$ws.Range("A21").Value = "EXTRA"
$ws.Range("A21").Font.Bold = $true

$ws.Range("B21").Value = "nn"
$ws.Range("C21").Value = "mo"
$ws.Range("D21").Value = "ft"
$ws.Range("E21").Value = "rt"
$ws.Range("F21").Value = "gs"
$ws.Range("G21").Value = "face"
$ws.Range("H21").Value = "dice"
$ws.Range("I21").Value = "mace"
$ws.Range("J21").Value = "cchvae"
$ws.Range("K21").Value = "juice"
$ws.Range("L21").Value = "Ijuice"

# --- Data rows 22-35: same dataset list as rows 2-15, appended below the new header ---
$datasets = @("adult","kdd_census","german","dutch","bank","credit","compass","diabetes","student","oulad","law","heart","synthetic_athlete","synthetic_disease")

for ($i = 0; $i -lt $datasets.Length; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = $datasets[$i]

    # B:E -> "Ubuntu", F:L -> "Home" (mirrors the layout of the original table)
    $ws.Cells.Item($row, 2).Value = "Ubuntu"
    $ws.Cells.Item($row, 3).Value = "Ubuntu"
    $ws.Cells.Item($row, 4).Value = "Ubuntu"
    $ws.Cells.Item($row, 5).Value = "Ubuntu"
    $ws.Cells.Item($row, 6).Value = "Home"
    $ws.Cells.Item($row, 7).Value = "Home"
    $ws.Cells.Item($row, 8).Value = "Home"
    $ws.Cells.Item($row, 9).Value = "Home"
    $ws.Cells.Item($row, 10).Value = "Home"
    $ws.Cells.Item($row, 11).Value = "Home"
    $ws.Cells.Item($row, 12).Value = "Home"
}

# First five datasets (adult, kdd_census, german, dutch, bank) already have the "mace"
# column (I) run & marked complete -> green fill, same as the "Full runs completed" flag
# used elsewhere in the sheet.
$ws.Range("I22:I26").Interior.Color = 5287936

# --- Selection / view bookkeeping ---
$ws.Range("J12").Select() | Out-Null

# --- Page layout touch that shows up in the saved file ---
$ws.PageSetup.Orientation = 1
